$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.764.22"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.632.15"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D5").Value = "215.25"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").Value = "0.0633"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "19.51"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.857.33"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "1.638.92"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D17").Value = "63.08"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "25.769.15"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "192.37"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "9.92"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").Value = "142.95"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "6.87"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").Value = "3.23"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").Value = "1.133.98"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").Value = "2.51"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").Value = "0.543"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "100.68"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "1.766.75"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").Value = "55.40"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").Value = "0.418"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  -7.02%  "
